$wb = $excel.ActiveWorkbook

# Update the batch size value on the Config sheet (B8: 500 -> 200)
$configSheet = $wb.Worksheets.Item("Config")
$configSheet.Range("B8").Value = 200

# Activate the Config sheet and select D18 (matches the final saved selection)
$configSheet.Activate()
$configSheet.Range("D18").Select()
